$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels to reflect generalized "container" terminology
$ws.Range("A2").Value = "Container Name:"
$ws.Range("B2").Value = "Taged Fish in Container"
$ws.Range("C2").Value = "# of Fish in Container"
$ws.Range("D2").Value = "Collection"

# Move the active selection from E2 to D3
$ws.Range("D3").Select()
